# outputs-HGR-r202/test-o__Haloplasmatales_A_split_pruned.xlsx
# "updated outputs-r202, previous copy of ful-path.csv"
#
# The header/label column (A1:C1 + A2) gets re-stamped with a (new) text
# number format -- same visual "@" text style the column already had, just
# written again as part of the refreshed export -- and the score in B2 is
# replaced with the freshly computed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text format to the header row and the row label, producing a
# distinct style entry for these cells (mirrors the diff's style index bump
# for A1, B1, C1 and A2).
$ws.Range("A1").NumberFormat = "@"
$ws.Range("B1").NumberFormat = "@"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"

# Updated prediction score for B2.
$ws.Range("B2").Value = 1684.2214090878099
